# Apply DPE_EPICP_P2 unit-conversion / algorithm corrections
# (Sodium, LDL, HDL, CHOL, TG unit conversions + AGE_CANCER algorithm fix)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 20: TG (Triglycerides) ---
$ws.Range("G20").Value = "operation"
$ws.Range("H20").Value = "corr_trigly/88.57"
$ws.Range("K20").Value = "compatible"

# --- Row 21: CHOL (Total cholesterol) ---
$ws.Range("G21").Value = "operation"
$ws.Range("H21").Value = "corr_chol/38.67"
$ws.Range("K21").Value = "compatible"

# --- Row 22: LDL cholesterol ---
$ws.Range("H22").Value = "(corr_chol - corr_hdl)/38.67"

# --- Row 23: HDL cholesterol ---
$ws.Range("G23").Value = "operation"
$ws.Range("H23").Value = "corr_hdl/38.67"
$ws.Range("K23").Value = "compatible"

# --- Row 60: AGE_CANCER ---
$ws.Range("F60").Value = "age0;d_recrui;dd_inccanc;dcens_canc;inccanc"
$ws.Range("G60").Value = "case_when"
$ws.Range("H60").Value = "case_when(inccanc == 1 ~ age0 + ((dd_inccanc - d_recrui)/365.25),`r`n                                           inccanc == 0 ~ age0 + ((dcens_canc - d_recrui)/365.25))"
$ws.Range("J60").Value = "partial"
$ws.Range("K60").Value = "proximate"

# --- Row 95: SODIUM ---
$ws.Range("G95").Value = "operation"
$ws.Range("H95").Value = "MNA*1000"
$ws.Range("K95").Value = "compatible"
